# Notifications.xlsx edit:
#  - Row 9 (TestCase with Jira/Description columns) had its Jira-id list and
#    Description list each lose one paired entry (the "comment on post"
#    item / OPQA-1397), and the row height shrinks to match the now-shorter
#    wrapped text.
#  - The active selection moves from C8 to C13.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newJiraIds = "OPQA-877||OPQA-1013||OPQA-215||OPQA-1395"
$newDescription = "Verify that user receives a notification when someone he is following  publishes a post||Verify that user is receiving notification when someone liked his post(aggregated notification)||Verify that user able to recevies a notification when other user commented on his post||Verify that all users receive notification when other user published a post and validate notification."

$ws.Range("B9").Value = $newJiraIds
$ws.Range("C9").Value = $newDescription

$ws.Rows.Item(9).RowHeight = 45

[void]$ws.Range("C13").Select()
